$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9119
$ws.Range("E2").Value = 159
$ws.Range("F2").Value = 159
$ws.Range("G2").Value = 258
$ws.Range("H2").Value = 193
$ws.Range("I2").Value = 204
$ws.Range("J2").Value = -12
$ws.Range("K2").Value = 8878
$ws.Range("L2").Value = 4381
$ws.Range("M2").Value = 4497
$ws.Range("N2").Value = 4492
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 357
$ws.Range("Q2").Value = 17
$ws.Range("R2").Value = -600
$ws.Range("S2").Value = 670
$ws.Range("T2").Value = 674
$ws.Range("U2").Value = -657
$ws.Range("V2").Value = 3182
$ws.Range("W2").Value = 1.74
$ws.Range("X2").Value = 2.11
$ws.Range("Y2").Value = 4.62
$ws.Range("Z2").Value = 2.34
$ws.Range("AA2").Value = 97.42
$ws.Range("AB2").Value = 1160.03
$ws.Range("AC2").Value = 2863
$ws.Range("AD2").Value = 39.82
$ws.Range("AE2").Value = 62912
$ws.Range("AF2").Value = 1.81
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 0.53
$ws.Range("AI2").Value = 20.96
$ws.Range("AJ2").Value = 7140000

# Row 3
$ws.Range("D3").Value = 10052
$ws.Range("E3").Value = 199
$ws.Range("F3").Value = 199
$ws.Range("G3").Value = 255
$ws.Range("H3").Value = 210
$ws.Range("I3").Value = 211
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 9142
$ws.Range("L3").Value = 4494
$ws.Range("M3").Value = 4648
$ws.Range("N3").Value = 4641
$ws.Range("O3").Value = 7
$ws.Range("P3").Value = 357
$ws.Range("Q3").Value = 479
$ws.Range("R3").Value = -656
$ws.Range("S3").Value = 190
$ws.Range("T3").Value = 610
$ws.Range("U3").Value = -132
$ws.Range("V3").Value = 3461
$ws.Range("W3").Value = 1.98
$ws.Range("X3").Value = 2.09
$ws.Range("Y3").Value = 4.63
$ws.Range("Z3").Value = 2.33
$ws.Range("AA3").Value = 96.69
$ws.Range("AB3").Value = 1202.05
$ws.Range("AC3").Value = 2958
$ws.Range("AD3").Value = 34.66
$ws.Range("AE3").Value = 64998
$ws.Range("AF3").Value = 1.58
$ws.Range("AG3").Value = 600
$ws.Range("AH3").Value = 0.59
$ws.Range("AI3").Value = 20.29
$ws.Range("AJ3").Value = 7140000

# Row 4
$ws.Range("D4").Value = 10211
$ws.Range("E4").Value = 270
$ws.Range("F4").Value = 270
$ws.Range("G4").Value = 255
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 174
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 9662
$ws.Range("L4").Value = 4877
$ws.Range("M4").Value = 4785
$ws.Range("N4").Value = 4767
$ws.Range("O4").Value = 18
$ws.Range("P4").Value = 357
$ws.Range("Q4").Value = 535
$ws.Range("R4").Value = -326
$ws.Range("S4").Value = -145
$ws.Range("T4").Value = 376
$ws.Range("U4").Value = 160
$ws.Range("V4").Value = 3540
$ws.Range("W4").Value = 2.65
$ws.Range("X4").Value = 1.71
$ws.Range("Y4").Value = 3.69
$ws.Range("Z4").Value = 1.86
$ws.Range("AA4").Value = 101.91
$ws.Range("AB4").Value = 1236.16
$ws.Range("AC4").Value = 2432
$ws.Range("AD4").Value = 27.22
$ws.Range("AE4").Value = 66769
$ws.Range("AF4").Value = 0.99
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 0.91
$ws.Range("AI4").Value = 24.67
$ws.Range("AJ4").Value = 7140000

# Row 5
$ws.Range("D5").Value = 11025
$ws.Range("E5").Value = 254
$ws.Range("F5").Value = 254
$ws.Range("G5").Value = 318
$ws.Range("H5").Value = 241
$ws.Range("I5").Value = 241
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 9509
$ws.Range("L5").Value = 4547
$ws.Range("M5").Value = 4962
$ws.Range("N5").Value = 4943
$ws.Range("O5").Value = 18
$ws.Range("P5").Value = 357
$ws.Range("Q5").Value = 139
$ws.Range("R5").Value = -312
$ws.Range("S5").Value = 37
$ws.Range("T5").Value = 328
$ws.Range("U5").Value = -189
$ws.Range("V5").Value = 3456
$ws.Range("W5").Value = 2.31
$ws.Range("X5").Value = 2.19
$ws.Range("Y5").Value = 4.96
$ws.Range("Z5").Value = 2.52
$ws.Range("AA5").Value = 91.63
$ws.Range("AB5").Value = 1284.1
$ws.Range("AC5").Value = 3371
$ws.Range("AD5").Value = 21.98
$ws.Range("AE5").Value = 69236
$ws.Range("AF5").Value = 1.07
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 0.8100000000000001
$ws.Range("AI5").Value = 17.8
$ws.Range("AJ5").Value = 7140000

# Row 6
$ws.Range("D6").Value = 12626
$ws.Range("E6").Value = 555
$ws.Range("F6").Value = 555
$ws.Range("G6").Value = 589
$ws.Range("H6").Value = 577
$ws.Range("I6").Value = 575
$ws.Range("K6").Value = 9825
$ws.Range("L6").Value = 4376
$ws.Range("M6").Value = 5448
$ws.Range("N6").Value = 5428
$ws.Range("P6").Value = 357
$ws.Range("Q6").Value = 1067
$ws.Range("R6").Value = -421
$ws.Range("S6").Value = -553
$ws.Range("T6").Value = 265
$ws.Range("U6").Value = 802
$ws.Range("V6").Value = 3050
$ws.Range("W6").Value = 4.4
$ws.Range("X6").Value = 4.57
$ws.Range("Y6").Value = 11.09
$ws.Range("Z6").Value = 5.97
$ws.Range("AA6").Value = 80.33
$ws.Range("AB6").Value = 1424.78
$ws.Range("AC6").Value = 8052
$ws.Range("AD6").Value = 24.16
$ws.Range("AE6").Value = 76024
$ws.Range("AF6").Value = 2.56
$ws.Range("AG6").Value = 850
$ws.Range("AH6").Value = 0.44
$ws.Range("AI6").Value = 10.56
$ws.Range("AJ6").Value = 7140000

# Row 7
$ws.Range("D7").Value = 14216
$ws.Range("E7").Value = 882
$ws.Range("G7").Value = 1042
$ws.Range("H7").Value = 793
$ws.Range("I7").Value = 795
$ws.Range("K7").Value = 11268
$ws.Range("L7").Value = 5140
$ws.Range("M7").Value = 6127
$ws.Range("N7").Value = 6119
$ws.Range("P7").Value = 359
$ws.Range("Q7").Value = 961
$ws.Range("R7").Value = -338
$ws.Range("S7").Value = -149
$ws.Range("T7").Value = 277
$ws.Range("U7").Value = 586
$ws.Range("W7").Value = 6.21
$ws.Range("X7").Value = 5.58
$ws.Range("Y7").Value = 13.77
$ws.Range("Z7").Value = 7.52
$ws.Range("AA7").Value = 83.88
$ws.Range("AC7").Value = 11135
$ws.Range("AD7").Value = 18.41
$ws.Range("AE7").Value = 85705
$ws.Range("AF7").Value = 2.39
$ws.Range("AG7").Value = 896
$ws.Range("AH7").Value = 0.44
$ws.Range("AI7").Value = 8.050000000000001

# Row 8
$ws.Range("D8").Value = 15459
$ws.Range("E8").Value = 1102
$ws.Range("G8").Value = 1258
$ws.Range("H8").Value = 963
$ws.Range("I8").Value = 960
$ws.Range("K8").Value = 12166
$ws.Range("L8").Value = 5152
$ws.Range("M8").Value = 7014
$ws.Range("N8").Value = 7012
$ws.Range("P8").Value = 359
$ws.Range("Q8").Value = 1001
$ws.Range("R8").Value = -247
$ws.Range("S8").Value = -278
$ws.Range("T8").Value = 200
$ws.Range("U8").Value = 828
$ws.Range("W8").Value = 7.13
$ws.Range("X8").Value = 6.23
$ws.Range("Y8").Value = 14.63
$ws.Range("Z8").Value = 8.220000000000001
$ws.Range("AA8").Value = 73.45999999999999
$ws.Range("AC8").Value = 13451
$ws.Range("AD8").Value = 15.24
$ws.Range("AE8").Value = 98211
$ws.Range("AF8").Value = 2.09
$ws.Range("AG8").Value = 939
$ws.Range("AH8").Value = 0.46
$ws.Range("AI8").Value = 6.98

# Row 9
$ws.Range("D9").Value = 16614
$ws.Range("E9").Value = 1325
$ws.Range("G9").Value = 1519
$ws.Range("H9").Value = 1164
$ws.Range("I9").Value = 1160
$ws.Range("K9").Value = 13282
$ws.Range("L9").Value = 5178
$ws.Range("M9").Value = 8105
$ws.Range("N9").Value = 8108
$ws.Range("P9").Value = 359
$ws.Range("Q9").Value = 1210
$ws.Range("R9").Value = -263
$ws.Range("S9").Value = -243
$ws.Range("T9").Value = 225
$ws.Range("U9").Value = 994
$ws.Range("W9").Value = 7.97
$ws.Range("X9").Value = 7.01
$ws.Range("Y9").Value = 15.35
$ws.Range("Z9").Value = 9.15
$ws.Range("AA9").Value = 63.89
$ws.Range("AC9").Value = 16249
$ws.Range("AD9").Value = 12.62
$ws.Range("AE9").Value = 113561
$ws.Range("AF9").Value = 1.81
$ws.Range("AG9").Value = 992
$ws.Range("AH9").Value = 0.48
$ws.Range("AI9").Value = 6.11
